# Generate Report for Handoff
# Updates the localization-status report: the b.md file has now been
# handed off again (new xlf files generated) while the handback for it
# is still the stale one -> flag "Ready for handoff" status + an
# out-of-date-handback Error Detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb1826a5a779122265f7a89a1ca83b73d614b95f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66fb31df5bc6a3b6c12f49eb976f2c23932b4307/e2e/b.md."

# --- Overview sheet: row 3 is the b.md entry ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 08:36:52"

# --- zh-cn sheet: row 3 is the b.md entry ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-28 08:36:47"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 is the b.md entry ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-28 08:36:52"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
